$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "保險" (insurance) -- 8th worksheet
# ----------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(8)

# Row 1 (header labels)
$ws8.Cells.Item(1,2).Value = "company"
$ws8.Cells.Item(1,3).Value = "name"
$ws8.Cells.Item(1,4).Value = "owner"
$ws8.Cells.Item(1,5).Value = "property_category"
$ws8.Cells.Item(1,6).Value = "category"
$ws8.Cells.Item(1,7).Value = "date"
$ws8.Cells.Item(1,8).Value = "legislator_name"
$ws8.Cells.Item(1,9).Value = "legislator_id"
$ws8.Cells.Item(1,10).Value = "source_file"
$ws8.Cells.Item(1,11).Value = "index"

# Row 2
$ws8.Cells.Item(2,2).Value = "三商美邦人壽"
$ws8.Cells.Item(2,3).Value = "20年繳費祥安終身壽險"
$ws8.Cells.Item(2,4).Value = "黃玉廷"
$ws8.Cells.Item(2,5).Value = "insurance"
$ws8.Cells.Item(2,6).Value = "normal"
$ws8.Cells.Item(2,7).NumberFormat = "@"
$ws8.Cells.Item(2,7).Value = "2012-04-24"
$ws8.Cells.Item(2,8).Value = "蔡其昌"
$ws8.Cells.Item(2,9).Value = 1377
$ws8.Cells.Item(2,10).Value = "tmp61ee1"
$ws8.Cells.Item(2,11).Value = 118

# Row 3
$ws8.Cells.Item(3,2).Value = "三商美邦人壽"
$ws8.Cells.Item(3,3).Value = "20年繳費祥安終身壽險"
$ws8.Cells.Item(3,4).Value = "黃玉廷"
$ws8.Cells.Item(3,5).Value = "insurance"
$ws8.Cells.Item(3,6).Value = "normal"
$ws8.Cells.Item(3,7).NumberFormat = "@"
$ws8.Cells.Item(3,7).Value = "2012-04-24"
$ws8.Cells.Item(3,8).Value = "蔡其昌"
$ws8.Cells.Item(3,9).Value = 1377
$ws8.Cells.Item(3,10).Value = "tmp61ee1"
$ws8.Cells.Item(3,11).Value = 119

# Row 4
$ws8.Cells.Item(4,2).Value = "國泰人壽"
$ws8.Cells.Item(4,3).Value = "添美盛美元終身壽險"
$ws8.Cells.Item(4,4).Value = "黃玉廷"
$ws8.Cells.Item(4,5).Value = "insurance"
$ws8.Cells.Item(4,6).Value = "normal"
$ws8.Cells.Item(4,7).NumberFormat = "@"
$ws8.Cells.Item(4,7).Value = "2012-04-24"
$ws8.Cells.Item(4,8).Value = "蔡其昌"
$ws8.Cells.Item(4,9).Value = 1377
$ws8.Cells.Item(4,10).Value = "tmp61ee1"
$ws8.Cells.Item(4,11).Value = 120

# Re-apply the shared header / data-row formatting (bold+border style used
# by B1:E1, plain style used by B2:E4) onto the newly populated cells so
# the new columns visually match the existing ones and no stray
# number-format survives on the date cells.
$ws8.Cells.Item(1,2).Copy()
$ws8.Range($ws8.Cells.Item(1,6), $ws8.Cells.Item(1,11)).PasteSpecial(-4122)

$ws8.Cells.Item(2,2).Copy()
$ws8.Range($ws8.Cells.Item(2,6), $ws8.Cells.Item(4,11)).PasteSpecial(-4122)

# ----------------------------------------------------------------------
# Sheet "債務" (debt) -- 9th worksheet
# ----------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item(9)

# Row 1 (header labels)
$ws9.Cells.Item(1,2).Value = "species"
$ws9.Cells.Item(1,3).Value = "debtor"
$ws9.Cells.Item(1,4).Value = "owner"
$ws9.Cells.Item(1,5).Value = "total"
$ws9.Cells.Item(1,6).Value = "register_date"
$ws9.Cells.Item(1,7).Value = "register_reason"
$ws9.Cells.Item(1,8).Value = "property_category"
$ws9.Cells.Item(1,9).Value = "category"
$ws9.Cells.Item(1,10).Value = "date"
$ws9.Cells.Item(1,11).Value = "legislator_name"
$ws9.Cells.Item(1,12).Value = "legislator_id"
$ws9.Cells.Item(1,13).Value = "source_file"
$ws9.Cells.Item(1,14).Value = "index"

# Row 2
$ws9.Cells.Item(2,2).Value = "房屋貸款"
$ws9.Cells.Item(2,3).Value = "蔡其昌"
$ws9.Cells.Item(2,4).Value = "安泰銀行豐原分行臺中市豐原區信義街"
$ws9.Cells.Item(2,5).Value = 4559066
$ws9.Cells.Item(2,6).Value = "91年04月04日"
$ws9.Cells.Item(2,7).Value = "親友與本人投資資金需求"
$ws9.Cells.Item(2,8).Value = "debt"
$ws9.Cells.Item(2,9).Value = "normal"
$ws9.Cells.Item(2,10).NumberFormat = "@"
$ws9.Cells.Item(2,10).Value = "2012-04-24"
$ws9.Cells.Item(2,11).Value = "蔡其昌"
$ws9.Cells.Item(2,12).Value = 1377
$ws9.Cells.Item(2,13).Value = "tmp61ee1"
$ws9.Cells.Item(2,14).Value = 130

# Row 3
$ws9.Cells.Item(3,2).Value = "房屋貸款"
$ws9.Cells.Item(3,3).Value = "黃玉廷"
$ws9.Cells.Item(3,4).Value = "第一銀行大甲分行臺中市大甲區順天路"
$ws9.Cells.Item(3,5).Value = 4926564
$ws9.Cells.Item(3,6).Value = "100年06月27日"
$ws9.Cells.Item(3,7).Value = "原房貸轉貸"
$ws9.Cells.Item(3,8).Value = "debt"
$ws9.Cells.Item(3,9).Value = "normal"
$ws9.Cells.Item(3,10).NumberFormat = "@"
$ws9.Cells.Item(3,10).Value = "2012-04-24"
$ws9.Cells.Item(3,11).Value = "蔡其昌"
$ws9.Cells.Item(3,12).Value = 1377
$ws9.Cells.Item(3,13).Value = "tmp61ee1"
$ws9.Cells.Item(3,14).Value = 131

# Row 4
$ws9.Cells.Item(4,2).Value = "房屋貸款"
$ws9.Cells.Item(4,3).Value = "黃玉廷"
$ws9.Cells.Item(4,4).Value = "國泰世華銀行沙鹿分行臺中市沙鹿區成功東街"
$ws9.Cells.Item(4,5).Value = 2349350
$ws9.Cells.Item(4,6).Value = "100年06月23日"
$ws9.Cells.Item(4,7).Value = "購置房屋"
$ws9.Cells.Item(4,8).Value = "debt"
$ws9.Cells.Item(4,9).Value = "normal"
$ws9.Cells.Item(4,10).NumberFormat = "@"
$ws9.Cells.Item(4,10).Value = "2012-04-24"
$ws9.Cells.Item(4,11).Value = "蔡其昌"
$ws9.Cells.Item(4,12).Value = 1377
$ws9.Cells.Item(4,13).Value = "tmp61ee1"
$ws9.Cells.Item(4,14).Value = 132

# Re-apply the shared header / data-row formatting onto the newly
# populated cells
$ws9.Cells.Item(1,2).Copy()
$ws9.Range($ws9.Cells.Item(1,8), $ws9.Cells.Item(1,14)).PasteSpecial(-4122)

$ws9.Cells.Item(2,2).Copy()
$ws9.Range($ws9.Cells.Item(2,8), $ws9.Cells.Item(4,14)).PasteSpecial(-4122)
